$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E2 value from 1 to 4
$ws.Range("E2").Value = 4

# Update selection from D2 to D9
$ws.Range("D9").Select()
